$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 ---
$ws.Range("A9").Value = 130654930
$ws.Range("B9").Value = 57884
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = 'Tretåig hackspett'
$ws.Range("G9").Value = 'Picoides tridactylus'
$ws.Range("H9").Value = '(Linnaeus, 1758)'
$ws.Range("J9").ClearContents()
$ws.Range("L9").Value = ''
$ws.Range("M9").Value = 'färska spår'
$ws.Range("Q9").Value = 440168
$ws.Range("R9").Value = 7053746
$ws.Range("AC9").Value = 'Ringhack (savhack), enstaka färska, några meter upp på en gran vid en hyggeskant.'
$ws.Range("AF9").ClearContents()
$ws.Range("AM9").Value = 'Trädstam på levande träd'
$ws.Range("AO9").Value = 'Stem on living tree # Picea abies'

# --- Row 10 ---
$ws.Range("A10").Value = 130654517
$ws.Range("B10").Value = 57884
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = 'Tretåig hackspett'
$ws.Range("G10").Value = 'Picoides tridactylus'
$ws.Range("H10").Value = '(Linnaeus, 1758)'
$ws.Range("J10").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("Q10").Value = 440178
$ws.Range("R10").Value = 7053979
$ws.Range("AC10").Value = 'Ringhack äldre'
$ws.Range("AF10").ClearContents()
$ws.Range("AH10").ClearContents()
$ws.Range("AJ10").ClearContents()
$ws.Range("AK10").ClearContents()
$ws.Range("AO10").ClearContents()
$ws.Range("AW10").Value = 'Benny Öwre'
$ws.Range("AX10").Value = 'Benny Öwre'

# --- Row 11 ---
$ws.Range("A11").Value = 130654938
$ws.Range("B11").Value = 79243
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = 'Garnlav'
$ws.Range("G11").Value = 'Alectoria sarmentosa'
$ws.Range("H11").Value = '(Ach.) Ach.'
$ws.Range("J11").Value = ''
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("Q11").Value = 440117
$ws.Range("R11").Value = 7053967
$ws.Range("AC11").ClearContents()
$ws.Range("AF11").Value = ''
$ws.Range("AM11").Value = 'Gren på levande träd'
$ws.Range("AO11").Value = 'Branch on living tree # Picea abies'

# --- Row 12 ---
$ws.Range("A12").Value = 130654941
$ws.Range("B12").Value = 79243
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = 'Garnlav'
$ws.Range("G12").Value = 'Alectoria sarmentosa'
$ws.Range("H12").Value = '(Ach.) Ach.'
$ws.Range("J12").Value = ''
$ws.Range("K12").Value = ''
$ws.Range("N12").Value = ''
$ws.Range("Q12").Value = 440134
$ws.Range("R12").Value = 7053783
$ws.Range("AC12").Value = 'Växer på en gran vid en hyggeskant.'
$ws.Range("AF12").Value = ''
$ws.Range("AH12").Value = 'Granskog'
$ws.Range("AJ12").Value = 'gran'
$ws.Range("AK12").Value = 'Picea abies'
$ws.Range("AO12").Value = 'Picea abies'
$ws.Range("AW12").Value = 'Kristian Zackrisson'
$ws.Range("AX12").Value = 'Kristian Zackrisson'

# --- Row 17 ---
$ws.Range("A17").Value = 130654935
$ws.Range("B17").Value = 79243
$ws.Range("E17").Value = 6425
$ws.Range("F17").Value = 'Garnlav'
$ws.Range("G17").Value = 'Alectoria sarmentosa'
$ws.Range("H17").Value = '(Ach.) Ach.'
$ws.Range("J17").Value = ''
$ws.Range("K17").Value = ''
$ws.Range("N17").Value = ''
$ws.Range("Q17").Value = 439862
$ws.Range("R17").Value = 7054226
$ws.Range("AC17").Value = 'Enstaka bålar på gran.'
$ws.Range("AF17").Value = ''
$ws.Range("AH17").Value = 'Granskog'
$ws.Range("AJ17").Value = 'gran'
$ws.Range("AK17").Value = 'Picea abies'
$ws.Range("AM17").Value = 'Gren på levande träd'
$ws.Range("AO17").Value = 'Branch on living tree # Picea abies'
$ws.Range("AW17").Value = 'Kristian Zackrisson'
$ws.Range("AX17").Value = 'Kristian Zackrisson'

# --- Row 18 ---
$ws.Range("A18").Value = 130654518
$ws.Range("B18").Value = 57884
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = 'Tretåig hackspett'
$ws.Range("G18").Value = 'Picoides tridactylus'
$ws.Range("H18").Value = '(Linnaeus, 1758)'
$ws.Range("J18").ClearContents()
$ws.Range("K18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("Q18").Value = 440177
$ws.Range("R18").Value = 7054022
$ws.Range("AC18").Value = 'Ringhack äldre'
$ws.Range("AF18").ClearContents()
$ws.Range("AH18").ClearContents()
$ws.Range("AJ18").ClearContents()
$ws.Range("AK18").ClearContents()
$ws.Range("AM18").ClearContents()
$ws.Range("AO18").ClearContents()
$ws.Range("AW18").Value = 'Benny Öwre'
$ws.Range("AX18").Value = 'Benny Öwre'
